$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New collection rows (MCH335-1 .. MCH335-4) appended below the header row.
# Columns: A identifier | C title | D date_s | E levelOfDescription
#          F extentAndMedium | G notes | H (blank, but still formatted)
# ---------------------------------------------------------------------------
$rows = @(
    @{ id = "MCH335-1"; title = "SOUTH- CITY OF CAPE TOWN: UTILITIES & WORK COMMITTEE AGENDAS FOR 1994-04-11, AMENNETIES & HEALTH COMMITTEE 1994-04-11, ITEMS LIST"; date = "1994" },
    @{ id = "MCH335-2"; title = "SOUTH- CITY OF CAPE TOWN: EXECUTIVE COMMITTEE, PLANNING COMMITTEE, ATHLONE & DISTRICT MANAGING COMMITTEE, WOODSTOCK/ WALMER WSTATE/ SALT RIVER"; date = $null },
    @{ id = "MCH335-3"; title = "SOUTH- CITY OF CAPE TOWN: EXECUTIVE COMMITTEE 1994-04-09, CAPE METRO NEGOTIATION FORUM 1994-03-17, RYLANDS ESTATE MANAGEMENT COMMITTEE 1994-03-15"; date = "1994" },
    @{ id = "MCH335-4"; title = "SOUTH- CITY OF CAPE TOWN: TOWN PLANNING, MANAGEMENT COMMITEE"; date = $null }
)

$levelOfDescription = "Series"
$extentAndMedium = "1 Box"
$notes = "LOCATION: 33G | GRAP COUNT NUMER: NONE"

$r = 2
foreach ($rowData in $rows) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellC = $ws.Cells.Item($r, 3)
    $cellD = $ws.Cells.Item($r, 4)
    $cellE = $ws.Cells.Item($r, 5)
    $cellF = $ws.Cells.Item($r, 6)
    $cellG = $ws.Cells.Item($r, 7)
    $cellH = $ws.Cells.Item($r, 8)

    $cellA.Value = $rowData.id
    $cellC.Value = $rowData.title
    if ($rowData.date) {
        # Force text storage (field is typed as text in the source data, "1994" not a number)
        $cellD.Value = "'" + $rowData.date
    }
    $cellE.Value = $levelOfDescription
    $cellF.Value = $extentAndMedium
    $cellG.Value = $notes

    foreach ($cell in @($cellA, $cellC, $cellD, $cellE, $cellF, $cellG, $cellH)) {
        $cell.Font.Name = "Calibri"
        $cell.Font.Size = 10
        $cell.Font.ThemeColor = 1
    }

    $ws.Rows.Item($r).RowHeight = 15.75

    $r = $r + 1
}

$ws.Rows.Item(1).RowHeight = 15.75

# ---------------------------------------------------------------------------
# Restore the frozen header pane + active selection (G14) as in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G14").Select()
